$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 984.3333
$ws.Range("I18").Value = 984.3333
$ws.Range("K18").Value = 984.3333
$ws.Range("M18").Value = -700.3333

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("I64").Value = 7082.6665
$ws.Range("J64").Value = 8923.666999999999
$ws.Range("K64").Value = 7082.6665
$ws.Range("L64").Value = 8923.666999999999
$ws.Range("M64").Value = -6834.6665
$ws.Range("N64").Value = -9419.666999999999

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("I67").Value = 7082.6665
$ws.Range("J67").Value = 8923.666999999999
$ws.Range("K67").Value = 7082.6665
$ws.Range("L67").Value = 8923.666999999999
$ws.Range("M67").Value = -6224.6665
$ws.Range("N67").Value = -10639.667

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H76").Value = 4499.4165
$ws.Range("I76").Value = 3032
$ws.Range("J76").Value = 5966.8335
$ws.Range("K76").Value = 3032
$ws.Range("L76").Value = 5966.8335
$ws.Range("M76").Value = -2717
$ws.Range("N76").Value = -6596.8335

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H79").Value = 4499.4165
$ws.Range("I79").Value = 3032
$ws.Range("J79").Value = 5966.8335
$ws.Range("K79").Value = 3032
$ws.Range("L79").Value = 5966.8335
$ws.Range("M79").Value = -1940
$ws.Range("N79").Value = -8150.8335

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 1701.2894
$ws.Range("I98").Value = 1670.8055
$ws.Range("J98").Value = 2250
$ws.Range("K98").Value = 1670.8055
$ws.Range("L98").Value = 2250
$ws.Range("M98").Value = -172.8054999999999
$ws.Range("N98").Value = -5246

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H106").Value = 37848.11
$ws.Range("I106").Value = 58986.4
$ws.Range("K106").Value = 58986.4
$ws.Range("M106").Value = -58355.4

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H116").Value = 2418.76
$ws.Range("I116").Value = 2613.5
$ws.Range("J116").Value = 2239
$ws.Range("K116").Value = 2613.5
$ws.Range("L116").Value = 2239
$ws.Range("M116").Value = 828.5
$ws.Range("N116").Value = -9123

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H122").Value = 1701.2894
$ws.Range("I122").Value = 1670.8055
$ws.Range("J122").Value = 2250
$ws.Range("K122").Value = 5012.416499999999
$ws.Range("L122").Value = 6750
$ws.Range("M122").Value = -2562.416499999999
$ws.Range("N122").Value = -11650

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1912.8462
$ws.Range("I32").Value = 1624.527
$ws.Range("K32").Value = 1624.527
$ws.Range("M32").Value = -1337.527

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2400.9697
$ws.Range("I61").Value = 1858.5264
$ws.Range("J61").Value = 3137.1428
$ws.Range("K61").Value = 1858.5264
$ws.Range("L61").Value = 3137.1428
$ws.Range("M61").Value = -1646.5264
$ws.Range("N61").Value = -3561.1428

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 8915.754999999999
$ws.Range("I132").Value = 4714.413
$ws.Range("J132").Value = 73336.336
$ws.Range("K132").Value = 14143.239
$ws.Range("L132").Value = 220009.008
$ws.Range("M132").Value = -11613.239
$ws.Range("N132").Value = -225069.008

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 2400.9697
$ws.Range("I136").Value = 1858.5264
$ws.Range("J136").Value = 3137.1428
$ws.Range("K136").Value = 5575.5792
$ws.Range("L136").Value = 9411.428400000001
$ws.Range("M136").Value = -3025.5792
$ws.Range("N136").Value = -14511.4284

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 3225.8572
$ws.Range("I86").Value = 2794.4614
$ws.Range("J86").Value = 3599.7334
$ws.Range("K86").Value = 2794.4614
$ws.Range("L86").Value = 3599.7334
$ws.Range("M86").Value = -1671.4614
$ws.Range("N86").Value = -5845.7334

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H89").Value = 3225.8572
$ws.Range("I89").Value = 2794.4614
$ws.Range("J89").Value = 3599.7334
$ws.Range("K89").Value = 13972.307
$ws.Range("L89").Value = 17998.667
$ws.Range("M89").Value = -8356.307000000001
$ws.Range("N89").Value = -29230.667

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 2019.9231
$ws.Range("I105").Value = 1975.8889
$ws.Range("J105").Value = 2119
$ws.Range("K105").Value = 1975.8889
$ws.Range("L105").Value = 2119
$ws.Range("M105").Value = -228.8888999999999
$ws.Range("N105").Value = -5613

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 2531.8867
$ws.Range("J134").Value = 8267.571
$ws.Range("L134").Value = 24802.713
$ws.Range("N134").Value = -29872.713

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H107").Value = 753.23914
$ws.Range("I107").Value = 336.8125
$ws.Range("K107").Value = 336.8125
$ws.Range("M107").Value = 1583.1875

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 12345868
$ws.Range("I2").Value = 140.86667
$ws.Range("J2").Value = 27778028
$ws.Range("K2").Value = 845.20002
$ws.Range("L2").Value = 166668168
$ws.Range("M2").Value = -732.20002
$ws.Range("N2").Value = -166668394

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1457.6
$ws.Range("I5").Value = 1219.9
$ws.Range("K5").Value = 3659.7
$ws.Range("M5").Value = -3547.7

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H135").Value = 1457.6
$ws.Range("I135").Value = 1219.9
$ws.Range("K135").Value = 10979.1
$ws.Range("M135").Value = -8444.1

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 48118.824
$ws.Range("I70").Value = 79180.875
$ws.Range("K70").Value = 79180.875
$ws.Range("M70").Value = -78910.875

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H73").Value = 48118.824
$ws.Range("I73").Value = 79180.875
$ws.Range("K73").Value = 79180.875
$ws.Range("M73").Value = -78244.875

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 664.75
$ws.Range("I46").Value = 724.5
$ws.Range("K46").Value = 724.5
$ws.Range("M46").Value = -536.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H128").Value = 62429
$ws.Range("J128").Value = 62429
$ws.Range("L128").Value = 62429
$ws.Range("N128").Value = -72389

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 3514.78
$ws.Range("I132").Value = 3479.3013
$ws.Range("J132").Value = 3688
$ws.Range("K132").Value = 10437.9039
$ws.Range("L132").Value = 11064
$ws.Range("M132").Value = -7907.903900000001
$ws.Range("N132").Value = -16124

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 4722.7
$ws.Range("J136").Value = 3728.0908
$ws.Range("L136").Value = 11184.2724
$ws.Range("N136").Value = -16284.2724

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 612
$ws.Range("I107").Value = 349.14285
$ws.Range("K107").Value = 1047.42855
$ws.Range("M107").Value = 872.5714499999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 3954.739
$ws.Range("I122").Value = 3675.7144
$ws.Range("K122").Value = 11027.1432
$ws.Range("M122").Value = -8577.143199999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 4028.889
$ws.Range("J136").Value = 3970.2856
$ws.Range("L136").Value = 11910.8568
$ws.Range("N136").Value = -17010.8568
